# Insert a new weekly observation row at row 217 (pushing the existing
# rows 217:322 down to 218:323) and populate it with the new data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(217).Insert()

$ws.Range("A217").Value = 4
$ws.Range("B217").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C217").Value = "Los Lagos"
$ws.Range("D217").Value = 44813
$ws.Range("E217").Value = 10
$ws.Range("F217").Value = 100112043
$ws.Range("G217").Value = "Pepino ensalada"
$ws.Range("H217").Value = "Sin especificar"
$ws.Range("I217").Value = "Primera"
$ws.Range("J217").Value = 400
$ws.Range("K217").Value = 30000
$ws.Range("L217").Value = 30000
$ws.Range("M217").Value = 30000
$ws.Range("N217").Value = "$/caja 60 unidades"
$ws.Range("O217").Value = "Región de Arica y Parinacota"
$ws.Range("P217").Value = 500
$ws.Range("Q217").Value = 60
$ws.Range("R217").Value = "Hortaliza"
